$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the value of D8 (keep formatting/style)
$ws.Range("D8").ClearContents()

# Update the selection to D8
$ws.Range("D8").Select()
